# Applies the "localdb" command-type addition to the #system reference sheet.
#
# Summary of the change:
#  1. A brand new column is inserted immediately before column N on the
#     (hidden) "#system" worksheet. It becomes the "localdb" column, holding
#     the names of 6 new SQL-related commands. Everything that used to live
#     in columns N..AC shifts one column to the right (O..AD).
#  2. The "target" list in column A (an alphabetically sorted catalogue of
#     every command-category name) gets a new entry, "localdb", inserted in
#     its correct alphabetical slot (row 14, just before "macro"). Everything
#     below shifts down one row - but only within column A, so the shift is
#     performed manually, cell by cell, rather than with a row-wide Insert.
#  3. The "web" list in column Y (formerly X, after the column insert above)
#     gets two new command names inserted in their correct alphabetical
#     slots: "scrollElement(...)" before "scrollLeft(...)" and
#     "scrollPage(...)" before "scrollRight(...)" - again shifting only
#     column Y.
#  4. All of the workbook-level defined names that point into the shifted
#     columns are corrected to point at their new locations, and a brand new
#     "localdb" defined name is added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Insert the new "localdb" column before column N, then populate it.
#    (A whole-column insert is safe here - it is meant to shift every
#    column from N rightwards.)
# ---------------------------------------------------------------------
$ws.Columns("N:N").Insert()

$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------
# 2. Insert "localdb" into the alphabetical "target" list (column A),
#    shifting A14:A29 down to A15:A30 - column A only. Since this
#    runtime's Range.Insert shifts the *entire* row (every column), the
#    shift is instead performed manually, one cell at a time, working
#    from the bottom up so no value is overwritten before being read.
# ---------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $src = $ws.Cells.Item($r, 1)
    $dst = $ws.Cells.Item($r + 1, 1)
    $dst.Value = $src.Value2
}
$ws.Range("A14").Value = "localdb"

# ---------------------------------------------------------------------
# 3. Insert the two new "web" function names into column Y (the old
#    column X, now shifted right because of the column insert above),
#    shifting only column Y's cells down, one insertion point at a time.
# ---------------------------------------------------------------------
for ($r = 125; $r -ge 101; $r--) {
    $src = $ws.Cells.Item($r, 25)
    $dst = $ws.Cells.Item($r + 1, 25)
    $dst.Value = $src.Value2
}
$ws.Range("Y101").Value = "scrollElement(locator,xOffset,yOffset)"

for ($r = 126; $r -ge 103; $r--) {
    $src = $ws.Cells.Item($r, 25)
    $dst = $ws.Cells.Item($r + 1, 25)
    $dst.Value = $src.Value2
}
$ws.Range("Y103").Value = "scrollPage(xOffset,yOffset)"

# ---------------------------------------------------------------------
# 4. Fix up the workbook-level defined names. Inserting the column does
#    not automatically repoint the existing defined names in this
#    runtime, so update them explicitly, then add "localdb".
# ---------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
